# Update market-price / profit columns (H:N) on each job sheet
# with refreshed Universalis pricing data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64: Forged from the Void
$ws.Range("H64").Value = 7500
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496

# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 7500
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716

# Row 113: Amaro Kart
$ws.Range("H113").Value = 4152.5
$ws.Range("I113").Value = 4005
$ws.Range("J113").Value = 4300
$ws.Range("K113").Value = 4005
$ws.Range("L113").Value = 4300
$ws.Range("M113").Value = -751
$ws.Range("N113").Value = -10808

# Row 121: Mindful Medicine
$ws.Range("H121").Value = 1500
$ws.Range("J121").Value = 1500
$ws.Range("L121").Value = 4500
$ws.Range("N121").Value = -7994

$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate
$ws.Range("H4").Value = 974.5
$ws.Range("I4").Value = 974.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 974.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -858.5
$ws.Range("N4").Value = ""

# Row 5: The Alloyed Truth
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""

# Row 39: Aurochs Star
$ws.Range("H39").Value = 2346.6
$ws.Range("I39").Value = 1429
$ws.Range("J39").Value = 6017
$ws.Range("K39").Value = 1429
$ws.Range("L39").Value = 6017
$ws.Range("M39").Value = -909
$ws.Range("N39").Value = -7057

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3555
$ws.Range("J61").Value = 3555
$ws.Range("L61").Value = 3555
$ws.Range("N61").Value = -3979

# Row 109: A Head of Demand
$ws.Range("H109").Value = 42999.332
$ws.Range("J109").Value = 42999.332
$ws.Range("L109").Value = 42999.332
$ws.Range("N109").Value = -45773.332

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 3555
$ws.Range("J136").Value = 3555
$ws.Range("L136").Value = 10665
$ws.Range("N136").Value = -15765

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = ""

# Row 38: The Naked Blade
$ws.Range("H38").Value = 31000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 31000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 31000
$ws.Range("M38").Value = ""
$ws.Range("N38").Value = -31832

$ws = $wb.Worksheets.Item("CRP")
# Row 99: O Pine
$ws.Range("H99").Value = 2257334.5
$ws.Range("I99").Value = 879854.75
$ws.Range("K99").Value = 879854.75
$ws.Range("M99").Value = -878356.75

# Row 126: A Better Conductor
$ws.Range("H126").Value = 2257334.5
$ws.Range("I126").Value = 879854.75
$ws.Range("K126").Value = 2639564.25
$ws.Range("M126").Value = -2637094.25

$ws = $wb.Worksheets.Item("CUL")
# Row 11: Putting the Squeeze On
$ws.Range("H11").Value = 250
$ws.Range("I11").Value = 250
$ws.Range("K11").Value = 750
$ws.Range("M11").Value = -610

# Row 14: Keep Your Powder Dry
$ws.Range("H14").Value = 1825.125
$ws.Range("I14").Value = 1825.125
$ws.Range("K14").Value = 5475.375
$ws.Range("M14").Value = -5302.375

# Row 25: Flakes for Friends
$ws.Range("H25").Value = 12
$ws.Range("I25").Value = 12
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 36
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 133
$ws.Range("N25").Value = ""

# Row 30: Picnic Panic
$ws.Range("H30").Value = 12
$ws.Range("I30").Value = 12
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 36
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 66
$ws.Range("N30").Value = ""

# Row 34: Fever Pitch
$ws.Range("H34").Value = 3819.25
$ws.Range("I34").Value = 2555
$ws.Range("K34").Value = 7665
$ws.Range("M34").Value = -7581

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers
$ws.Range("H2").Value = 165.66667
$ws.Range("I2").Value = 194.58333
$ws.Range("J2").Value = 107.833336
$ws.Range("K2").Value = 194.58333
$ws.Range("L2").Value = 107.833336
$ws.Range("M2").Value = -81.58332999999999
$ws.Range("N2").Value = -333.833336

# Row 11: A Ringing Success
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = ""

# Row 59: Sew Not Doing This
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -21166

# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 100023384
$ws.Range("I70").Value = 38669.332
$ws.Range("J70").Value = 250000450
$ws.Range("K70").Value = 38669.332
$ws.Range("L70").Value = 250000450
$ws.Range("M70").Value = -38399.332
$ws.Range("N70").Value = -250000990

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 100023384
$ws.Range("I73").Value = 38669.332
$ws.Range("J73").Value = 250000450
$ws.Range("K73").Value = 38669.332
$ws.Range("L73").Value = 250000450
$ws.Range("M73").Value = -37733.332
$ws.Range("N73").Value = -250002322

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 2179
$ws.Range("I126").Value = 2258.7144
$ws.Range("J126").Value = 1900
$ws.Range("K126").Value = 6776.1432
$ws.Range("L126").Value = 5700
$ws.Range("M126").Value = -4306.1432
$ws.Range("N126").Value = -10640

$ws = $wb.Worksheets.Item("LTW")
# Row 12: A Place to Call Helm
$ws.Range("H12").Value = 1664.6666
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1664.6666
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1664.6666
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = -2004.6666

# Row 42: Slave to Fashion
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = ""

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = ""

# Row 49: First They Came for the Heretics
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""

# Row 110: Breeches of Trust
$ws.Range("H110").Value = 31190.334
$ws.Range("J110").Value = 31190.334
$ws.Range("L110").Value = 31190.334
$ws.Range("N110").Value = -39370.334

$ws = $wb.Worksheets.Item("WVR")
# Row 103: To the Tops
$ws.Range("H103").Value = 16854.75
$ws.Range("J103").Value = 16854.75
$ws.Range("L103").Value = 16854.75
$ws.Range("N103").Value = -19198.75

# Row 107: Flax Wax
$ws.Range("H107").Value = 483.7
$ws.Range("I107").Value = 299.8
$ws.Range("J107").Value = 667.6
$ws.Range("K107").Value = 899.4000000000001
$ws.Range("L107").Value = 2002.8
$ws.Range("M107").Value = 1020.6
$ws.Range("N107").Value = -5842.8
